$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.213.35"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.913.26"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "484.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.32"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000346"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "4.534.20"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "3.906.65"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.23"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.25"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "68.260.92"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.52"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.04"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.67"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +18.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.71"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("E27").Value = "  +11.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "715.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("E34").Value = "  +15.91%  "
$ws.Range("E35").Value = "  +6.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "41.69"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "61.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.403"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +20.85%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +17.25%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -4.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0489"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E51").Value = "  +27.33%  "
